$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Vtn"
$ws.Cells.Item(2,3).Value = "Plaur"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 12.75206033333333
$ws.Cells.Item(2,8).Value = 38.256181
$ws.Cells.Item(2,9).Value = 0.1573122343381959
$ws.Cells.Item(2,10).Value = 0.157312234338196
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 27.81717833333333
$ws.Cells.Item(2,14).Value = 83.451535
$ws.Cells.Item(2,15).Value = 0.4044740580248731
$ws.Cells.Item(2,16).Value = 0.4044740580248732
$ws.Cells.Item(2,17).Value = 354.7263364097594
$ws.Cells.Item(2,18).Value = 3192.537027687835
$ws.Cells.Item(2,19).Value = 0.06362871779972991
$ws.Cells.Item(2,20).Value = 0.06362871779972992

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Vtn"
$ws.Cells.Item(3,3).Value = "Plaur"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 12.75206033333333
$ws.Cells.Item(3,8).Value = 38.256181
$ws.Cells.Item(3,9).Value = 0.1573122343381959
$ws.Cells.Item(3,10).Value = 0.157312234338196
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 4.827410666666666
$ws.Cells.Item(3,14).Value = 14.482232
$ws.Cells.Item(3,15).Value = 0.0701926830500802
$ws.Cells.Item(3,16).Value = 0.0701926830500802
$ws.Cells.Item(3,17).Value = 61.55943207511022
$ws.Cells.Item(3,18).Value = 554.034888675992
$ws.Cells.Item(3,19).Value = 0.01104216780480093
$ws.Cells.Item(3,20).Value = 0.01104216780480093

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Vtn"
$ws.Cells.Item(4,3).Value = "Plaur"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 12.75206033333333
$ws.Cells.Item(4,8).Value = 38.256181
$ws.Cells.Item(4,9).Value = 0.1573122343381959
$ws.Cells.Item(4,10).Value = 0.157312234338196
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 32.06242466666667
$ws.Cells.Item(4,14).Value = 96.187274
$ws.Cells.Item(4,15).Value = 0.4662018145637509
$ws.Cells.Item(4,16).Value = 0.466201814563751
$ws.Cells.Item(4,17).Value = 408.8619737822882
$ws.Cells.Item(4,18).Value = 3679.757764040594
$ws.Cells.Item(4,19).Value = 0.07333924910154495
$ws.Cells.Item(4,20).Value = 0.07333924910154498

$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Vtn"
$ws.Cells.Item(5,3).Value = "Plaur"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 12.75206033333333
$ws.Cells.Item(5,8).Value = 38.256181
$ws.Cells.Item(5,9).Value = 0.1573122343381959
$ws.Cells.Item(5,10).Value = 0.157312234338196
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 4.066688333333334
$ws.Cells.Item(5,14).Value = 12.200065
$ws.Cells.Item(5,15).Value = 0.05913144436129575
$ws.Cells.Item(5,16).Value = 0.05913144436129575
$ws.Cells.Item(5,17).Value = 51.85865498352945
$ws.Cells.Item(5,18).Value = 466.7278948517651
$ws.Cells.Item(5,19).Value = 0.009302099632120153
$ws.Cells.Item(5,20).Value = 0.009302099632120153

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Vtn"
$ws.Cells.Item(6,3).Value = "Plaur"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 20.35396833333334
$ws.Cells.Item(6,8).Value = 61.06190500000001
$ws.Cells.Item(6,9).Value = 0.2510910513649196
$ws.Cells.Item(6,10).Value = 0.2510910513649196
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 27.81717833333333
$ws.Cells.Item(6,14).Value = 83.451535
$ws.Cells.Item(6,15).Value = 0.4044740580248731
$ws.Cells.Item(6,16).Value = 0.4044740580248732
$ws.Cells.Item(6,17).Value = 566.189966919353
$ws.Cells.Item(6,18).Value = 5095.709702274176
$ws.Cells.Item(6,19).Value = 0.1015598164793009
$ws.Cells.Item(6,20).Value = 0.1015598164793009

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Vtn"
$ws.Cells.Item(7,3).Value = "Plaur"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 20.35396833333334
$ws.Cells.Item(7,8).Value = 61.06190500000001
$ws.Cells.Item(7,9).Value = 0.2510910513649196
$ws.Cells.Item(7,10).Value = 0.2510910513649196
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 4.827410666666666
$ws.Cells.Item(7,14).Value = 14.482232
$ws.Cells.Item(7,15).Value = 0.0701926830500802
$ws.Cells.Item(7,16).Value = 0.0701926830500802
$ws.Cells.Item(7,17).Value = 98.2569638413289
$ws.Cells.Item(7,18).Value = 884.3126745719601
$ws.Cells.Item(7,19).Value = 0.01762475458516921
$ws.Cells.Item(7,20).Value = 0.01762475458516921

$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Vtn"
$ws.Cells.Item(8,3).Value = "Plaur"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 20.35396833333334
$ws.Cells.Item(8,8).Value = 61.06190500000001
$ws.Cells.Item(8,9).Value = 0.2510910513649196
$ws.Cells.Item(8,10).Value = 0.2510910513649196
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 32.06242466666667
$ws.Cells.Item(8,14).Value = 96.187274
$ws.Cells.Item(8,15).Value = 0.4662018145637509
$ws.Cells.Item(8,16).Value = 0.466201814563751
$ws.Cells.Item(8,17).Value = 652.597576355219
$ws.Cells.Item(8,18).Value = 5873.378187196971
$ws.Cells.Item(8,19).Value = 0.1170591037670455
$ws.Cells.Item(8,20).Value = 0.1170591037670455

$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Vtn"
$ws.Cells.Item(9,3).Value = "Plaur"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 20.35396833333334
$ws.Cells.Item(9,8).Value = 61.06190500000001
$ws.Cells.Item(9,9).Value = 0.2510910513649196
$ws.Cells.Item(9,10).Value = 0.2510910513649196
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 4.066688333333334
$ws.Cells.Item(9,14).Value = 12.200065
$ws.Cells.Item(9,15).Value = 0.05913144436129575
$ws.Cells.Item(9,16).Value = 0.05913144436129575
$ws.Cells.Item(9,17).Value = 82.77324555820282
$ws.Cells.Item(9,18).Value = 744.9592100238252
$ws.Cells.Item(9,19).Value = 0.014847376533404
$ws.Cells.Item(9,20).Value = 0.014847376533404

$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Vtn"
$ws.Cells.Item(10,3).Value = "Plaur"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 0.3333333333333333
$ws.Cells.Item(10,7).Value = 0.004706
$ws.Cells.Item(10,8).Value = 0.014118
$ws.Cells.Item(10,9).Value = 0.00005805425597465284
$ws.Cells.Item(10,10).Value = 0.00005805425597465285
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 27.81717833333333
$ws.Cells.Item(10,14).Value = 83.451535
$ws.Cells.Item(10,15).Value = 0.4044740580248731
$ws.Cells.Item(10,16).Value = 0.4044740580248732
$ws.Cells.Item(10,17).Value = 0.1309076412366667
$ws.Cells.Item(10,18).Value = 1.17816877113
$ws.Cells.Item(10,19).Value = 0.00002348144049968257
$ws.Cells.Item(10,20).Value = 0.00002348144049968258

$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Vtn"
$ws.Cells.Item(11,3).Value = "Plaur"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(11,6).Value = 0.3333333333333333
$ws.Cells.Item(11,7).Value = 0.004706
$ws.Cells.Item(11,8).Value = 0.014118
$ws.Cells.Item(11,9).Value = 0.00005805425597465284
$ws.Cells.Item(11,10).Value = 0.00005805425597465285
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 4.827410666666666
$ws.Cells.Item(11,14).Value = 14.482232
$ws.Cells.Item(11,15).Value = 0.0701926830500802
$ws.Cells.Item(11,16).Value = 0.0701926830500802
$ws.Cells.Item(11,17).Value = 0.02271779459733333
$ws.Cells.Item(11,18).Value = 0.204460151376
$ws.Cells.Item(11,19).Value = 0.000004074983989337031
$ws.Cells.Item(11,20).Value = 0.000004074983989337032

$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Vtn"
$ws.Cells.Item(12,3).Value = "Plaur"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 1
$ws.Cells.Item(12,6).Value = 0.3333333333333333
$ws.Cells.Item(12,7).Value = 0.004706
$ws.Cells.Item(12,8).Value = 0.014118
$ws.Cells.Item(12,9).Value = 0.00005805425597465284
$ws.Cells.Item(12,10).Value = 0.00005805425597465285
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 32.06242466666667
$ws.Cells.Item(12,14).Value = 96.187274
$ws.Cells.Item(12,15).Value = 0.4662018145637509
$ws.Cells.Item(12,16).Value = 0.466201814563751
$ws.Cells.Item(12,17).Value = 0.1508857704813333
$ws.Cells.Item(12,18).Value = 1.357971934332
$ws.Cells.Item(12,19).Value = 0.00002706499947853163
$ws.Cells.Item(12,20).Value = 0.00002706499947853164

$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Vtn"
$ws.Cells.Item(13,3).Value = "Plaur"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(13,6).Value = 0.3333333333333333
$ws.Cells.Item(13,7).Value = 0.004706
$ws.Cells.Item(13,8).Value = 0.014118
$ws.Cells.Item(13,9).Value = 0.00005805425597465284
$ws.Cells.Item(13,10).Value = 0.00005805425597465285
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 4.066688333333334
$ws.Cells.Item(13,14).Value = 12.200065
$ws.Cells.Item(13,15).Value = 0.05913144436129575
$ws.Cells.Item(13,16).Value = 0.05913144436129575
$ws.Cells.Item(13,17).Value = 0.01913783529666667
$ws.Cells.Item(13,18).Value = 0.17224051767
$ws.Cells.Item(13,19).Value = 0.000003432832007101606
$ws.Cells.Item(13,20).Value = 0.000003432832007101606

$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Vtn"
$ws.Cells.Item(14,3).Value = "Plaur"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 47.951367
$ws.Cells.Item(14,8).Value = 143.854101
$ws.Cells.Item(14,9).Value = 0.5915386600409097
$ws.Cells.Item(14,10).Value = 0.5915386600409098
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 27.81717833333333
$ws.Cells.Item(14,14).Value = 83.451535
$ws.Cells.Item(14,15).Value = 0.4044740580248731
$ws.Cells.Item(14,16).Value = 0.4044740580248732
$ws.Cells.Item(14,17).Value = 1333.871727166115
$ws.Cells.Item(14,18).Value = 12004.84554449504
$ws.Cells.Item(14,19).Value = 0.2392620423053426
$ws.Cells.Item(14,20).Value = 0.2392620423053427

$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Vtn"
$ws.Cells.Item(15,3).Value = "Plaur"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 47.951367
$ws.Cells.Item(15,8).Value = 143.854101
$ws.Cells.Item(15,9).Value = 0.5915386600409097
$ws.Cells.Item(15,10).Value = 0.5915386600409098
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 4.827410666666666
$ws.Cells.Item(15,14).Value = 14.482232
$ws.Cells.Item(15,15).Value = 0.0701926830500802
$ws.Cells.Item(15,16).Value = 0.0701926830500802
$ws.Cells.Item(15,17).Value = 231.480940537048
$ws.Cells.Item(15,18).Value = 2083.328464833432
$ws.Cells.Item(15,19).Value = 0.04152168567612072
$ws.Cells.Item(15,20).Value = 0.04152168567612072

$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Vtn"
$ws.Cells.Item(16,3).Value = "Plaur"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 47.951367
$ws.Cells.Item(16,8).Value = 143.854101
$ws.Cells.Item(16,9).Value = 0.5915386600409097
$ws.Cells.Item(16,10).Value = 0.5915386600409098
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 32.06242466666667
$ws.Cells.Item(16,14).Value = 96.187274
$ws.Cells.Item(16,15).Value = 0.4662018145637509
$ws.Cells.Item(16,16).Value = 0.466201814563751
$ws.Cells.Item(16,17).Value = 1537.437092101186
$ws.Cells.Item(16,18).Value = 13836.93382891068
$ws.Cells.Item(16,19).Value = 0.2757763966956819
$ws.Cells.Item(16,20).Value = 0.275776396695682

$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Vtn"
$ws.Cells.Item(17,3).Value = "Plaur"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 47.951367
$ws.Cells.Item(17,8).Value = 143.854101
$ws.Cells.Item(17,9).Value = 0.5915386600409097
$ws.Cells.Item(17,10).Value = 0.5915386600409098
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 4.066688333333334
$ws.Cells.Item(17,14).Value = 12.200065
$ws.Cells.Item(17,15).Value = 0.05913144436129575
$ws.Cells.Item(17,16).Value = 0.05913144436129575
$ws.Cells.Item(17,17).Value = 195.0032647462851
$ws.Cells.Item(17,18).Value = 1755.029382716566
$ws.Cells.Item(17,19).Value = 0.0349785353637645
$ws.Cells.Item(17,20).Value = 0.0349785353637645

